# Add 2022-Q3 data
# 1) Insert a new worksheet "2022-Q3" right after "总计" (becomes the 2nd tab).
# 2) Populate it with the new quarter's fund holdings table.
# 3) Insert a corresponding new top row into the "总计" summary sheet and
#    push the existing rows down by one.

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet right after "总计"
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $totalSheet)
$newSheet.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 2) Fill in the new sheet's data
# ---------------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $newSheet.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$rows = @(
    @("012463", "博时成长优势混合A",        "15.96", "84.33", "4.12", "0.6576", 5),
    @("004823", "上投摩根安裕回报混合A",      "5.35",  "25.71", "3.47", "0.1856", 1),
    @("004824", "上投摩根安裕回报混合C",      "4.91",  "25.71", "3.47", "0.1704", 1),
    @("011034", "南方宝恒混合C",             "14.78", "25.32", "0.81", "0.1197", 6),
    @("011033", "南方宝恒混合A",             "14.52", "25.32", "0.81", "0.1176", 6),
    @("010742", "南方宁悦一年持有期混合A",     "12.41", "26.26", "0.80", "0.0993", 8),
    @("012464", "博时成长优势混合C",          "0.67",  "84.33", "4.12", "0.0276", 5),
    @("010743", "南方宁悦一年持有期混合C",     "1.57",  "26.26", "0.80", "0.0126", 8)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $newSheet.Cells.Item($r, 1).Value = $i

    # Columns B-G are stored as text (keeps leading zeros / exact decimals)
    $bg = $newSheet.Range("B" + $r + ":G" + $r)
    $bg.NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $data[0]
    $newSheet.Cells.Item($r, 3).Value = $data[1]
    $newSheet.Cells.Item($r, 4).Value = $data[2]
    $newSheet.Cells.Item($r, 5).Value = $data[3]
    $newSheet.Cells.Item($r, 6).Value = $data[4]
    $newSheet.Cells.Item($r, 7).Value = $data[5]

    # Column H is a real number
    $newSheet.Cells.Item($r, 8).Value = $data[6]
}

# Apply header/index-column look (bold, thin border, centered/top) matching
# the styling used on every other quarterly sheet. (Union ranges don't
# reliably format every area in this runtime, so style each area separately.)
foreach ($rng in @($newSheet.Range("B1:H1"), $newSheet.Range("A2:A9"))) {
    $rng.Font.Bold = $true
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
}

# ---------------------------------------------------------------------------
# 3) Update the "总计" sheet: insert the new quarter at the top and shift
#    the existing rows down by one.
# ---------------------------------------------------------------------------

# Extend formatting of the new last row (row 8) by copying row 7's format.
$totalSheet.Range("A7:D7").Copy($totalSheet.Range("A8:D8"))

# Shift existing data rows down (only the real data columns B:D; column A is
# just a running index that gets rewritten below).
for ($r = 7; $r -ge 2; $r--) {
    $totalSheet.Cells.Item($r + 1, 2).Value = $totalSheet.Cells.Item($r, 2).Value2
    $totalSheet.Cells.Item($r + 1, 3).Value = $totalSheet.Cells.Item($r, 3).Value2
    $totalSheet.Cells.Item($r + 1, 4).Value = $totalSheet.Cells.Item($r, 4).Value2
}

# Write the new quarter into the now-empty row 2.
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 8
$totalSheet.Cells.Item(2, 4).Value = 1.39

# Re-number the running index column (0-based).
for ($r = 2; $r -le 8; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

# Restore the original active sheet/selection.
$totalSheet.Activate()
$totalSheet.Range("A1").Select()
